$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.771852374076843
$ws.Range("B1").Value = 1.968438267707825
$ws.Range("C1").Value = 2.233006238937378
$ws.Range("D1").Value = 2.744485855102539
$ws.Range("E1").Value = 1.386294841766357
